$wb = $excel.ActiveWorkbook

# Rename the first sheet (F-SW-SD-06 -> S-SW-SC-06) and keep the Print_Area
# defined name pointing at the renamed sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "S-SW-SC-06"
$ws1.PageSetup.PrintArea = '$A$1:$G$31'

# Hide the helper "Sheet2" tab.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Visible = $false

# Scroll the visible sheet's window so row 24 is at the top (topLeftCell = A24),
# keeping the existing E24 selection.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
